$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A87").Value = 43987
